$d = $word.ActiveDocument

$replacements = @(
    @{old="23×52="; new="66×79="},
    @{old="51×65="; new="54×88="},
    @{old="21×50="; new="34×25="},
    @{old="36×44="; new="93×53="},
    @{old="40×28="; new="59×72="},
    @{old="20×71="; new="95×92="},
    @{old="23×76="; new="41×68="},
    @{old="44×49="; new="46×23="},
    @{old="44×75="; new="64×18="},
    @{old="80×68="; new="55×67="},
    @{old="51×45="; new="27×77="},
    @{old="65×44="; new="97×30="},
    @{old="87×48="; new="86×51="},
    @{old="15×65="; new="62×72="},
    @{old="14×75="; new="65×99="},
    @{old="98×16="; new="99×92="},
    @{old="97×26="; new="26×49="},
    @{old="18×39="; new="49×41="},
    @{old="56×40="; new="49×78="},
    @{old="65×80="; new="87×22="},
    @{old="86×73="; new="95×70="},
    @{old="97×77="; new="46×88="},
    @{old="64×59="; new="91×61="},
    @{old="44×56="; new="42×24="},
    @{old="86×64="; new="71×81="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
